$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "67.025.85"
$ws.Range("E2").Value = "  +0.19%  "
Set-TextCell $ws.Range("D3") "3.210.93"
$ws.Range("E3").Value = "  -0.59%  "
Set-TextCell $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell $ws.Range("D5") "576.80"
$ws.Range("E5").Value = "  -1.96%  "
Set-TextCell $ws.Range("D6") "141.91"
$ws.Range("E6").Value = "  -7.13%  "
Set-TextCell $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  +0.13%  "
Set-TextCell $ws.Range("D8") "3.205.88"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  -5.25%  "
Set-TextCell $ws.Range("D11") "6.25"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  -3.46%  "
Set-TextCell $ws.Range("D14") "35.73"
$ws.Range("E14").Value = "  -7.66%  "
Set-TextCell $ws.Range("D15") "3.718.68"
$ws.Range("E15").Value = "  -0.96%  "
Set-TextCell $ws.Range("D16") "66.929.91"
$ws.Range("E16").Value = "  +0.15%  "
Set-TextCell $ws.Range("D17") "3.211.83"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  -1.89%  "
Set-TextCell $ws.Range("D19") "6.83"
$ws.Range("E19").Value = "  -1.79%  "
Set-TextCell $ws.Range("D20") "500.31"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E23").Value = "  -5.45%  "
Set-TextCell $ws.Range("D24") "81.67"
$ws.Range("E24").Value = "  -3.27%  "
Set-TextCell $ws.Range("D25") "12.84"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("E28").Value = "  -5.80%  "
Set-TextCell $ws.Range("D29") "27.68"
$ws.Range("E29").Value = "  -3.46%  "
Set-TextCell $ws.Range("D30") "7.57"
$ws.Range("E30").Value = "  -6.07%  "
$ws.Range("E31").Value = "  +2.01%  "
Set-TextCell $ws.Range("D32") "2.52"
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("E33").Value = "  -0.12%  "
Set-TextCell $ws.Range("D34") "513.41"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D35") "54.42"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D36") "6.05"
$ws.Range("E36").Value = "  -5.64%  "
$ws.Range("E37").Value = "  -7.96%  "
Set-TextCell $ws.Range("D38") "0.0412"
$ws.Range("E38").Value = "  -3.28%  "
Set-TextCell $ws.Range("D39") "0.0811"
$ws.Range("E39").Value = "  -4.81%  "
Set-TextCell $ws.Range("D40") "8.57"
$ws.Range("E40").Value = "  -6.77%  "
$ws.Range("E41").Value = "  -6.03%  "
Set-TextCell $ws.Range("D42") "2.854.20"
$ws.Range("E42").Value = "  -1.53%  "
Set-TextCell $ws.Range("D43") "2.50"
$ws.Range("E43").Value = "  -10.63%  "
$ws.Range("E45").Value = "  -3.62%  "
Set-TextCell $ws.Range("D46") "121.59"
$ws.Range("E46").Value = "  +0.09%  "
Set-TextCell $ws.Range("D47") "24.74"
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("E48").Value = "  -6.86%  "
Set-TextCell $ws.Range("D49") "0.0₃0523"
$ws.Range("E49").Value = "  -10.85%  "
Set-TextCell $ws.Range("D50") "0.109"
$ws.Range("E50").Value = "  -2.83%  "
Set-TextCell $ws.Range("D51") "2.10"
$ws.Range("E51").Value = "  -12.58%  "
